$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) values for B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values for B2:E2
$ws.Range("B2").Value = 457.7521169658944
$ws.Range("C2").Value = 387.34516811258061
$ws.Range("D2").Value = 455.3124955684969
$ws.Range("E2").Value = 383.1533072526862

# Update row 3 (STR) values for B3:E3
$ws.Range("B3").Value = 460.74879719978543
$ws.Range("C3").Value = 387.17347803187351
$ws.Range("D3").Value = 456.28044220366763
$ws.Range("E3").Value = 392.6553446345344

# Update the selection to match B1:E3
$ws.Range("B1:E3").Select()
